$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.021.65'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.423.85'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.15'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '89.36'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('E7').Value = '  -2.81%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.497'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0835'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.67%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '31.91'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.77%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.109'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.795.33'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.72'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.85'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.404.65'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.772'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.44%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '40.957.77'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0925'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.94%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.26'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.31'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.03'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.07'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.69'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').Value = '  -3.03%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.58'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '34.88'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.08%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '155.69'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.56%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.27'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -2.74%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0746'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.55%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.00'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.97%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '16.69'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.59%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.114'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.79'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('E40').Value = '  -2.81%  '
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.995.98'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.27'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -8.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '18.87'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.70%  '
$ws.Range('E45').Value = '  -3.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.89'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.52'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.36%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.652.28'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '95.12'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '73.66'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '52.32'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.47%  '
